$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Unique occurrences - safe to use a simple Find & Replace across the whole document.
Replace-Text "70×87=6090" "78×30=2340"
Replace-Text "23×59=1357" "73×53=3869"
Replace-Text "92×62=5704" "99×20=1980"
Replace-Text "90×16=1440" "24×99=2376"
Replace-Text "49×62=3038" "25×32=800"
Replace-Text "93×73=6789" "55×55=3025"
Replace-Text "87×86=7482" "46×96=4416"
Replace-Text "22×82=1804" "38×69=2622"
Replace-Text "74×90=6660" "88×30=2640"
Replace-Text "86×59=5074" "67×24=1608"
Replace-Text "17×69=1173" "47×68=3196"
Replace-Text "12×50=600" "16×79=1264"
Replace-Text "54×73=3942" "60×15=900"
Replace-Text "23×81=1863" "66×29=1914"
Replace-Text "50×24=1200" "87×29=2523"
Replace-Text "77×16=1232" "36×76=2736"
Replace-Text "39×42=1638" "24×46=1104"
Replace-Text "19×54=1026" "18×82=1476"
Replace-Text "57×54=3078" "46×68=3128"
Replace-Text "70×46=3220" "65×71=4615"
Replace-Text "41×22=902" "60×17=1020"
Replace-Text "84×75=6300" "90×48=4320"
Replace-Text "47×74=3478" "40×53=2120"

# "74×74=5476" appears twice in the table (Row 10 Col 5, and Row 15 Col 3) and is
# replaced by two different values, so address each table cell explicitly
# rather than risk a global Find touching the wrong instance.
$table = $d.Tables(1)
$table.Cell(10, 5).Range.Text = "28×84=2352"
$table.Cell(15, 3).Range.Text = "53×86=4558"
